# Applies the "ifo gdp component analysis preprocessing" edit:
# For each of rows 45-53, update the last existing numeric value in the row
# and append one more quarter-ahead error value (Q9 column) to the right of it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell value, then write the new cell immediately to its right.
$ws.Range("I45").Value = 0.5011245128056051
$ws.Range("J45").Value = 0.2858677898194339

$ws.Range("H46").Value = 0.4852787037784192
$ws.Range("I46").Value = 0.2775335613519331

$ws.Range("G47").Value = 0.4539510573947921
$ws.Range("H47").Value = 0.2743085116504074

$ws.Range("F48").Value = 0.4663391832225094
$ws.Range("G48").Value = 0.2534447081011285

$ws.Range("E49").Value = 0.4814444548743619
$ws.Range("F49").Value = 0.2766837437271186

$ws.Range("D50").Value = 0.4184715358843989
$ws.Range("E50").Value = 0.2867219094086165

$ws.Range("C51").Value = 0.5177895860664353
$ws.Range("D51").Value = 0.1751453671933744

$ws.Range("B52").Value = 0.5618492773058843
$ws.Range("C52").Value = 0.1965658720679752

$ws.Range("B53").Value = 0.4328090033804217
